$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2481.0645
$ws.Range("I113").Value = 2566.6667
$ws.Range("J113").Value = 2460.52
$ws.Range("K113").Value = 2566.6667
$ws.Range("L113").Value = 2460.52
$ws.Range("M113").Value = 687.3332999999998
$ws.Range("N113").Value = -8968.52

$ws.Range("H116").Value = 5153449
$ws.Range("I116").Value = 5450663.5
$ws.Range("J116").Value = 1733.3334
$ws.Range("K116").Value = 5450663.5
$ws.Range("L116").Value = 1733.3334
$ws.Range("M116").Value = -5447221.5
$ws.Range("N116").Value = -8617.3334

$ws.Range("H132").Value = 6016.6875
$ws.Range("I132").Value = 6943.9473
$ws.Range("J132").Value = 4661.4614
$ws.Range("K132").Value = 20831.8419
$ws.Range("L132").Value = 13984.3842
$ws.Range("M132").Value = -18301.8419
$ws.Range("N132").Value = -19044.3842

$ws.Range("H133").Value = 43900
$ws.Range("J133").Value = 43900
$ws.Range("L133").Value = 43900
$ws.Range("N133").Value = -54020

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11015.333
$ws.Range("I32").Value = 2644.7646
$ws.Range("J32").Value = 31343.857
$ws.Range("K32").Value = 2644.7646
$ws.Range("L32").Value = 31343.857
$ws.Range("M32").Value = -2357.7646
$ws.Range("N32").Value = -31917.857

$ws.Range("H61").Value = 1200.7407
$ws.Range("I61").Value = 1132.8125
$ws.Range("K61").Value = 1132.8125
$ws.Range("M61").Value = -920.8125

$ws.Range("H74").Value = 4809649.5
$ws.Range("I74").Value = 5683403
$ws.Range("J74").Value = 4007
$ws.Range("K74").Value = 5683403
$ws.Range("L74").Value = 4007
$ws.Range("M74").Value = -5682529
$ws.Range("N74").Value = -5755

$ws.Range("H77").Value = 4809649.5
$ws.Range("I77").Value = 5683403
$ws.Range("J77").Value = 4007
$ws.Range("K77").Value = 28417015
$ws.Range("L77").Value = 20035
$ws.Range("M77").Value = -28412647
$ws.Range("N77").Value = -28771

$ws.Range("H132").Value = 1482.2174
$ws.Range("I132").Value = 1286.9412
$ws.Range("J132").Value = 2035.5
$ws.Range("K132").Value = 3860.8236
$ws.Range("L132").Value = 6106.5
$ws.Range("M132").Value = -1330.8236
$ws.Range("N132").Value = -11166.5

$ws.Range("H136").Value = 1200.7407
$ws.Range("I136").Value = 1132.8125
$ws.Range("K136").Value = 3398.4375
$ws.Range("M136").Value = -848.4375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 51.5
$ws.Range("I80").Value = 19
$ws.Range("J80").Value = 54
$ws.Range("K80").Value = 19
$ws.Range("L80").Value = 54
$ws.Range("M80").Value = 979
$ws.Range("N80").Value = -2050

$ws.Range("H83").Value = 51.5
$ws.Range("I83").Value = 19
$ws.Range("J83").Value = 54
$ws.Range("K83").Value = 95
$ws.Range("L83").Value = 270
$ws.Range("M83").Value = 4897
$ws.Range("N83").Value = -10254

$ws.Range("H99").Value = 1515.3636
$ws.Range("I99").Value = 1361.125
$ws.Range("J99").Value = 1926.6666
$ws.Range("K99").Value = 1361.125
$ws.Range("L99").Value = 1926.6666
$ws.Range("M99").Value = 136.875
$ws.Range("N99").Value = -4922.6666

$ws.Range("H134").Value = 2010.5186
$ws.Range("I134").Value = 1729.1305
$ws.Range("K134").Value = 5187.3915
$ws.Range("M134").Value = -2652.3915

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4253.706
$ws.Range("I99").Value = 6385.2856
$ws.Range("J99").Value = 2761.6
$ws.Range("K99").Value = 6385.2856
$ws.Range("L99").Value = 2761.6
$ws.Range("M99").Value = -4887.2856
$ws.Range("N99").Value = -5757.6

$ws.Range("H126").Value = 4253.706
$ws.Range("I126").Value = 6385.2856
$ws.Range("J126").Value = 2761.6
$ws.Range("K126").Value = 19155.8568
$ws.Range("L126").Value = 8284.799999999999
$ws.Range("M126").Value = -16685.8568
$ws.Range("N126").Value = -13224.8

$ws.Range("H132").Value = 1678.6316
$ws.Range("I132").Value = 1339.6
$ws.Range("J132").Value = 2950
$ws.Range("K132").Value = 4018.8
$ws.Range("L132").Value = 8850
$ws.Range("M132").Value = -1488.8
$ws.Range("N132").Value = -13910

$ws.Range("H134").Value = 3639.6365
$ws.Range("I134").Value = 4071.5
$ws.Range("K134").Value = 12214.5
$ws.Range("M134").Value = -9679.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 392.61703
$ws.Range("I5").Value = 349.14633
$ws.Range("J5").Value = 689.6667
$ws.Range("K5").Value = 1047.43899
$ws.Range("L5").Value = 2069.0001
$ws.Range("M5").Value = -935.4389899999999
$ws.Range("N5").Value = -2293.0001

$ws.Range("H12").Value = 34.35
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 40.235294
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 120.705882
$ws.Range("M12").Value = 170
$ws.Range("N12").Value = -466.705882

$ws.Range("H113").Value = 1016.85
$ws.Range("I113").Value = 488.54544
$ws.Range("J113").Value = 1662.5555
$ws.Range("K113").Value = 1465.63632
$ws.Range("L113").Value = 4987.666499999999
$ws.Range("M113").Value = 704.3636799999999
$ws.Range("N113").Value = -9327.666499999999

$ws.Range("H122").Value = 407.04544
$ws.Range("I122").Value = 386.375
$ws.Range("J122").Value = 418.85715
$ws.Range("K122").Value = 3477.375
$ws.Range("L122").Value = 3769.71435
$ws.Range("M122").Value = -1027.375
$ws.Range("N122").Value = -8669.71435

$ws.Range("H135").Value = 392.61703
$ws.Range("I135").Value = 349.14633
$ws.Range("J135").Value = 689.6667
$ws.Range("K135").Value = 3142.31697
$ws.Range("L135").Value = 6207.0003
$ws.Range("M135").Value = -607.3169699999999
$ws.Range("N135").Value = -11277.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1267.1852
$ws.Range("I102").Value = 1252.174
$ws.Range("J102").Value = 1353.5
$ws.Range("K102").Value = 1252.174
$ws.Range("L102").Value = 1353.5
$ws.Range("M102").Value = 369.826
$ws.Range("N102").Value = -4597.5

$ws.Range("H122").Value = 3046.9363
$ws.Range("I122").Value = 2556.7837
$ws.Range("J122").Value = 4860.5
$ws.Range("K122").Value = 7670.3511
$ws.Range("L122").Value = 14581.5
$ws.Range("M122").Value = -5220.3511
$ws.Range("N122").Value = -19481.5

$ws.Range("H132").Value = 2094.111
$ws.Range("I132").Value = 1748
$ws.Range("J132").Value = 2371
$ws.Range("K132").Value = 5244
$ws.Range("L132").Value = 7113
$ws.Range("M132").Value = -2714
$ws.Range("N132").Value = -12173

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 14888100
$ws.Range("I132").Value = 23157572
$ws.Range("J132").Value = 3050.2666
$ws.Range("K132").Value = 69472716
$ws.Range("L132").Value = 9150.799800000001
$ws.Range("M132").Value = -69470186
$ws.Range("N132").Value = -14210.7998
